# Automatische test-sync: 2025-07-23 21:26:50
# Append a new incoming-mail log row to the "Logs" sheet and bump the
# matching rollup counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$newRow = 10

$logs.Cells.Item($newRow, 1).Value = "Wat zijn jullie openingstijden?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Cells.Item($newRow, 4).Value = "Openingstijden / Locatie"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value = "2025-07-23 21:25:56"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Ja"

# Setting a value with embedded newlines auto-expands the row height in
# this runtime; AutoFit settles it back to the sheet's default (matching
# every other data row, which carries no explicit row height).
$logs.Rows.Item($newRow).AutoFit()

# Extend the existing conditional-formatting rules (Categorie/Beantwoord/
# Handmatig opvolgen/Automatisch afgehandeld/Hergebruikt antwoord columns)
# so they keep covering the data range now that it runs through row 10.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "9")
    $newRange = $logs.Range($col + "2:" + $col + $newRow)
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Bump the "Openingstijden / Locatie" rollup count on the Dashboard sheet.
$dash.Cells.Item(2, 2).Value = 9
